# Auto-generated Excel COM-interop script
# Applies numeric cell updates to the "Lamia_Profits" workbook (8 sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as captured by the scheduled-runner commit diff.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 5845.8335
$ws.Range("J17").Value = 6034.7827
$ws.Range("L17").Value = 18104.3481
$ws.Range("N17").Value = -18440.3481
# Row 106
$ws.Range("H106").Value = 2630.7942
$ws.Range("I106").Value = 1519.5217
$ws.Range("K106").Value = 1519.5217
$ws.Range("M106").Value = -888.5217
# Row 109
$ws.Range("H109").Value = 87625
$ws.Range("J109").Value = 87625
$ws.Range("L109").Value = 87625
$ws.Range("N109").Value = -90399
# Row 111
$ws.Range("H111").Value = 1454
$ws.Range("I111").Value = 1454
$ws.Range("K111").Value = 4362
$ws.Range("M111").Value = -1295
# Row 112
$ws.Range("H112").Value = 1819
$ws.Range("J112").Value = 1876.9
$ws.Range("L112").Value = 5630.700000000001
$ws.Range("N112").Value = -7846.700000000001
# Row 125
$ws.Range("H125").Value = 4374.375
$ws.Range("J125").Value = 4699.3335
$ws.Range("L125").Value = 42294.0015
$ws.Range("N125").Value = -47214.0015

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5391.3223
$ws.Range("I32").Value = 4078.6538
$ws.Range("K32").Value = 4078.6538
$ws.Range("M32").Value = -3791.6538
# Row 61
$ws.Range("H61").Value = 21002.6
$ws.Range("I61").Value = 18999.75
$ws.Range("J61").Value = 29014
$ws.Range("K61").Value = 18999.75
$ws.Range("L61").Value = 29014
$ws.Range("M61").Value = -18787.75
$ws.Range("N61").Value = -29438
# Row 74
$ws.Range("H74").Value = 13891635
$ws.Range("I74").Value = 22224320
$ws.Range("K74").Value = 22224320
$ws.Range("M74").Value = -22223446
# Row 77
$ws.Range("H77").Value = 13891635
$ws.Range("I77").Value = 22224320
$ws.Range("K77").Value = 111121600
$ws.Range("M77").Value = -111117232
# Row 97
$ws.Range("H97").Value = 3563
$ws.Range("I97").Value = 3563
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 3563
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -3067
$ws.Range("N97").ClearContents()
# Row 132
$ws.Range("H132").Value = 5772.909
$ws.Range("I132").Value = 6003.2812
$ws.Range("J132").Value = 5158.5835
$ws.Range("K132").Value = 18009.8436
$ws.Range("L132").Value = 15475.7505
$ws.Range("M132").Value = -15479.8436
$ws.Range("N132").Value = -20535.7505
# Row 136
$ws.Range("H136").Value = 21002.6
$ws.Range("I136").Value = 18999.75
$ws.Range("J136").Value = 29014
$ws.Range("K136").Value = 56999.25
$ws.Range("L136").Value = 87042
$ws.Range("M136").Value = -54449.25
$ws.Range("N136").Value = -92142

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1801.2858
$ws.Range("I94").Value = 1801.2858
$ws.Range("K94").Value = 1801.2858
$ws.Range("M94").Value = -1350.2858
# Row 96
$ws.Range("H96").Value = 15809.333
$ws.Range("I96").Value = 15809.333
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 15809.333
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -13063.333
$ws.Range("N96").ClearContents()
# Row 134
$ws.Range("H134").Value = 3262.75
$ws.Range("I134").Value = 2013.9
$ws.Range("K134").Value = 6041.700000000001
$ws.Range("M134").Value = -3506.700000000001

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 27664.047
$ws.Range("I31").Value = 1939.7391
$ws.Range("J31").Value = 58804
$ws.Range("K31").Value = 1939.7391
$ws.Range("L31").Value = 58804
$ws.Range("M31").Value = -1644.7391
$ws.Range("N31").Value = -59394
# Row 34
$ws.Range("H34").Value = 27664.047
$ws.Range("I34").Value = 1939.7391
$ws.Range("J34").Value = 58804
$ws.Range("K34").Value = 1939.7391
$ws.Range("L34").Value = 58804
$ws.Range("M34").Value = -1737.7391
$ws.Range("N34").Value = -59208
# Row 58
$ws.Range("H58").Value = 3296.5789
$ws.Range("I58").Value = 1713.9375
$ws.Range("J58").Value = 11737.333
$ws.Range("K58").Value = 1713.9375
$ws.Range("L58").Value = 11737.333
$ws.Range("M58").Value = -1510.9375
$ws.Range("N58").Value = -12143.333
# Row 99
$ws.Range("H99").Value = 2366
$ws.Range("J99").Value = 2366
$ws.Range("L99").Value = 2366
$ws.Range("N99").Value = -5362
# Row 126
$ws.Range("H126").Value = 2366
$ws.Range("J126").Value = 2366
$ws.Range("L126").Value = 7098
$ws.Range("N126").Value = -12038
# Row 132
$ws.Range("H132").Value = 3193.1667
$ws.Range("I132").Value = 2670.3462
$ws.Range("K132").Value = 8011.0386
$ws.Range("M132").Value = -5481.0386
# Row 136
$ws.Range("H136").Value = 3296.5789
$ws.Range("I136").Value = 1713.9375
$ws.Range("J136").Value = 11737.333
$ws.Range("K136").Value = 5141.8125
$ws.Range("L136").Value = 35211.999
$ws.Range("M136").Value = -2591.8125
$ws.Range("N136").Value = -40311.999

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 3300.1428
$ws.Range("J5").Value = 10002.5
$ws.Range("L5").Value = 30007.5
$ws.Range("N5").Value = -30231.5
# Row 37
$ws.Range("H37").Value = 330000.53
$ws.Range("J37").Value = 330000.53
$ws.Range("L37").Value = 990001.5900000001
$ws.Range("N37").Value = -990225.5900000001
# Row 60
$ws.Range("H60").Value = 37950310
$ws.Range("I60").Value = 55555630
$ws.Range("J60").Value = 2739666.8
$ws.Range("K60").Value = 166666890
$ws.Range("L60").Value = 8219000.399999999
$ws.Range("M60").Value = -166666639
$ws.Range("N60").Value = -8219502.399999999
# Row 121
$ws.Range("H121").Value = 1333.3077
$ws.Range("I121").Value = 1476.4286
$ws.Range("J121").Value = 1166.3334
$ws.Range("K121").Value = 4429.2858
$ws.Range("L121").Value = 3499.0002
$ws.Range("M121").Value = -3119.2858
$ws.Range("N121").Value = -6119.0002
# Row 135
$ws.Range("H135").Value = 3300.1428
$ws.Range("J135").Value = 10002.5
$ws.Range("L135").Value = 90022.5
$ws.Range("N135").Value = -95092.5
# Row 140
$ws.Range("H140").Value = 2293.4583
$ws.Range("I140").Value = 1476.05
$ws.Range("K140").Value = 4428.15
$ws.Range("M140").Value = 751.8500000000004

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 315333.5
$ws.Range("I80").Value = 387179.78
$ws.Range("J80").Value = 3999.6667
$ws.Range("K80").Value = 387179.78
$ws.Range("L80").Value = 3999.6667
$ws.Range("M80").Value = -386181.78
$ws.Range("N80").Value = -5995.6667
# Row 83
$ws.Range("H83").Value = 315333.5
$ws.Range("I83").Value = 387179.78
$ws.Range("J83").Value = 3999.6667
$ws.Range("K83").Value = 1935898.9
$ws.Range("L83").Value = 19998.3335
$ws.Range("M83").Value = -1930906.9
$ws.Range("N83").Value = -29982.3335
# Row 113
$ws.Range("H113").Value = 2935.8667
$ws.Range("I113").Value = 1800.25
$ws.Range("J113").Value = 4233.7144
$ws.Range("K113").Value = 1800.25
$ws.Range("L113").Value = 4233.7144
$ws.Range("M113").Value = 369.75
$ws.Range("N113").Value = -8573.714400000001
# Row 122
$ws.Range("H122").Value = 5176.434
$ws.Range("I122").Value = 3608.04
$ws.Range("J122").Value = 6576.7856
$ws.Range("K122").Value = 10824.12
$ws.Range("L122").Value = 19730.3568
$ws.Range("M122").Value = -8374.119999999999
$ws.Range("N122").Value = -24630.3568
# Row 132
$ws.Range("H132").Value = 3143.5186
$ws.Range("I132").Value = 2585.1177
$ws.Range("J132").Value = 4092.8
$ws.Range("K132").Value = 7755.353099999999
$ws.Range("L132").Value = 12278.4
$ws.Range("M132").Value = -5225.353099999999
$ws.Range("N132").Value = -17338.4
# Row 134
$ws.Range("H134").Value = 73695.2
$ws.Range("J134").Value = 73695.2
$ws.Range("L134").Value = 221085.6
$ws.Range("N134").Value = -226155.6
# Row 136
$ws.Range("H136").Value = 29074.066
$ws.Range("J136").Value = 27079.357
$ws.Range("L136").Value = 81238.071
$ws.Range("N136").Value = -86338.071

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 14573.833
$ws.Range("I2").Value = 20000
$ws.Range("J2").Value = 10698
$ws.Range("K2").Value = 20000
$ws.Range("L2").Value = 10698
$ws.Range("M2").Value = -19888
$ws.Range("N2").Value = -10922
# Row 22
$ws.Range("H22").Value = 2343.6765
$ws.Range("I22").Value = 2008.174
$ws.Range("J22").Value = 3045.182
$ws.Range("K22").Value = 2008.174
$ws.Range("L22").Value = 3045.182
$ws.Range("M22").Value = -1713.174
$ws.Range("N22").Value = -3635.182
# Row 27
$ws.Range("H27").Value = 2343.6765
$ws.Range("I27").Value = 2008.174
$ws.Range("J27").Value = 3045.182
$ws.Range("K27").Value = 2008.174
$ws.Range("L27").Value = 3045.182
$ws.Range("M27").Value = -1901.174
$ws.Range("N27").Value = -3259.182
# Row 40
$ws.Range("H40").Value = 10432.35
$ws.Range("I40").Value = 10565.1
$ws.Range("K40").Value = 10565.1
$ws.Range("M40").Value = -10429.1
# Row 82
$ws.Range("H82").Value = 6191.1816
$ws.Range("I82").Value = 1549.5714
$ws.Range("J82").Value = 14314
$ws.Range("K82").Value = 1549.5714
$ws.Range("L82").Value = 14314
$ws.Range("M82").Value = -1188.5714
$ws.Range("N82").Value = -15036
# Row 85
$ws.Range("H85").Value = 6191.1816
$ws.Range("I85").Value = 1549.5714
$ws.Range("J85").Value = 14314
$ws.Range("K85").Value = 1549.5714
$ws.Range("L85").Value = 14314
$ws.Range("M85").Value = -301.5714
$ws.Range("N85").Value = -16810
# Row 93
$ws.Range("H93").Value = 1903.7368
$ws.Range("I93").Value = 1872.4706
$ws.Range("J93").Value = 2169.5
$ws.Range("K93").Value = 1872.4706
$ws.Range("L93").Value = 2169.5
$ws.Range("M93").Value = -624.4706000000001
$ws.Range("N93").Value = -4665.5
# Row 122
$ws.Range("H122").Value = 369625.53
$ws.Range("J122").Value = 9944.25
$ws.Range("L122").Value = 29832.75
$ws.Range("N122").Value = -34732.75
# Row 128
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
# Row 132
$ws.Range("H132").Value = 4731.3213
$ws.Range("I132").Value = 4836.2095
$ws.Range("K132").Value = 14508.6285
$ws.Range("M132").Value = -11978.6285

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 86
$ws.Range("H86").Value = 48999.5
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 48999.5
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 48999.5
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -51245.5
# Row 89
$ws.Range("H89").Value = 48999.5
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 48999.5
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 244997.5
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -256229.5
# Row 122
$ws.Range("H122").Value = 1756.66
$ws.Range("I122").Value = 1270.6316
$ws.Range("K122").Value = 3811.8948
$ws.Range("M122").Value = -1361.8948
